$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.900.02'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.812.69'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.22'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4661'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3661'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07341'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8667'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.30'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.794.92'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.504'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.72'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008693'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.63'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.923.88'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.294'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.050.07'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.894'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.39'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.26'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.156'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.271'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.50'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08917'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7542'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.156'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.484'
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.911'
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  -1.92%  '
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.994'
$ws.Range('E39').Value = '  +2.37%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.205'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5295'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.285'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.389'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4857'
$ws.Range('E46').Value = '  -2.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.41'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.12'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.658'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('E51').Value = '  -0.05%  '
